$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 137, pushing the existing rows 137-171 down to 138-172.
$ws.Rows("137").Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A137").Value = 4
$ws.Range("B137").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C137").Value = "Los Lagos"
$ws.Range("D137").Value = 44508
$ws.Range("E137").Value = 10
$ws.Range("F137").Value = 100112044
$ws.Range("G137").Value = "Perejil"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 80
$ws.Range("K137").Value = 5000
$ws.Range("L137").Value = 5000
$ws.Range("M137").Value = 5000
$ws.Range("N137").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O137").Value = "Región Metropolitana"
$ws.Range("P137").Value = 1667
$ws.Range("Q137").Value = 3
$ws.Range("R137").Value = "Hortaliza"
